# "Generate Report for Handoff"
# - Status changes from "In Translation" to "Ready for handoff" (Overview!E2/F2,
#   zh-cn!C2, de-de!C2 all share this text)
# - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamp is
#   refreshed (Overview!G2 & de-de!H2 share 00:43:23 -> 00:43:59; zh-cn!H2 moves
#   00:43:18 -> 00:43:55)
# - The Status/zh-cn/de-de columns get a little wider to fit the new text

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Refreshed timestamps ---
$wsOverview.Range("G2").Value = "2016-09-04 00:43:59"
$wsDeDe.Range("H2").Value = "2016-09-04 00:43:59"
$wsZhCn.Range("H2").Value = "2016-09-04 00:43:55"

# --- Column width adjustments (Status / zh-cn / de-de columns widened) ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.33   # zh-cn column
$wsOverview.Columns.Item(6).ColumnWidth = 16.33   # de-de column
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33        # Status column
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33        # Status column
